$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.970.00'
$ws.Range('E2').Value = '  +0.24%  '
$ws.Range('D3').Value = '2.754.51'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.25'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.77%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '158.44'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.65%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.604'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.01%  '
$ws.Range('E9').Value = '  -1.97%  '
$ws.Range('E10').Value = '  +1.83%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.383'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.46%  '
$ws.Range('E12').Value = '  -16.98%  '
$ws.Range('D13').Value = '3.236.37'
$ws.Range('E13').Value = '  +0.46%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.90'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.30%  '
$ws.Range('D15').Value = '63.594.40'
$ws.Range('E15').Value = '  -0.27%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000152'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.91%  '
$ws.Range('D17').Value = '2.752.35'
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.19'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.98%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.86'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '357.39'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.75'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.11%  '
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.535'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.34'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.26%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.171'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.67'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.23%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.34%  '
$ws.Range('D28').Value = '0.0₃0915'
$ws.Range('E28').Value = '  -1.28%  '
$ws.Range('E29').Value = '  -3.80%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.22'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.38%  '
$ws.Range('E31').Value = '  -1.47%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '169.57'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.48%  '
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.95'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.86%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '20.23'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.99%  '
$ws.Range('E35').Value = '  +1.24%  '
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.81'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.28%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.28'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.79%  '
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '335.54'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.78%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.19'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.37%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.03'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.40'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.80'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0591'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0256'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.79%  '
$ws.Range('E47').Value = '  -0.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '135.50'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.37%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.626'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.42%  '
$ws.Range('E50').Value = '  +0.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.05'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.05%  '
